$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 5000
$ws.Range("J13").Value = 5000
$ws.Range("L13").Value = 5000
$ws.Range("N13").Value = -5338
$ws.Range("H16").Value = 4999.5
$ws.Range("I16").Value = 4999.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 4999.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -4769.5
$ws.Range("N16").ClearContents()
$ws.Range("H62").Value = 4117
$ws.Range("I62").Value = 3645.5
$ws.Range("K62").Value = 3645.5
$ws.Range("M62").Value = -3021.5
$ws.Range("H65").Value = 4117
$ws.Range("I65").Value = 3645.5
$ws.Range("K65").Value = 18227.5
$ws.Range("M65").Value = -15107.5
$ws.Range("H97").Value = 954.25
$ws.Range("J97").Value = 954.25
$ws.Range("L97").Value = 2862.75
$ws.Range("N97").Value = -3854.75
$ws.Range("H98").Value = 661.53845
$ws.Range("I98").Value = 679.75
$ws.Range("J98").Value = 443
$ws.Range("K98").Value = 679.75
$ws.Range("L98").Value = 443
$ws.Range("M98").Value = 818.25
$ws.Range("N98").Value = -3439
$ws.Range("H116").Value = 8340158
$ws.Range("I116").Value = 6750
$ws.Range("J116").Value = 16673566
$ws.Range("K116").Value = 6750
$ws.Range("L116").Value = 16673566
$ws.Range("M116").Value = -3308
$ws.Range("N116").Value = -16680450
$ws.Range("H122").Value = 661.53845
$ws.Range("I122").Value = 679.75
$ws.Range("J122").Value = 443
$ws.Range("K122").Value = 2039.25
$ws.Range("L122").Value = 1329
$ws.Range("M122").Value = 410.75
$ws.Range("N122").Value = -6229
$ws.Range("H127").Value = 923.5
$ws.Range("I127").Value = 1047.75
$ws.Range("J127").Value = 675
$ws.Range("K127").Value = 3143.25
$ws.Range("L127").Value = 2025
$ws.Range("M127").Value = 1816.75
$ws.Range("N127").Value = -11945
$ws.Range("H129").Value = 2526.45
$ws.Range("I129").Value = 1564.3125
$ws.Range("J129").Value = 6375
$ws.Range("K129").Value = 4692.9375
$ws.Range("L129").Value = 19125
$ws.Range("M129").Value = 307.0625
$ws.Range("N129").Value = -29125
$ws.Range("H132").Value = 1683.7241
$ws.Range("I132").Value = 1746.4445
$ws.Range("J132").Value = 837
$ws.Range("K132").Value = 5239.333500000001
$ws.Range("L132").Value = 2511
$ws.Range("M132").Value = -2709.333500000001
$ws.Range("N132").Value = -7571
$ws.Range("H138").Value = 1709.0193
$ws.Range("I138").Value = 1219.4688
$ws.Range("J138").Value = 2492.3
$ws.Range("K138").Value = 3658.4064
$ws.Range("L138").Value = 7476.900000000001
$ws.Range("M138").Value = 1481.5936
$ws.Range("N138").Value = -17756.9
$ws.Range("H141").Value = 2771.5715
$ws.Range("I141").Value = 2215.5386
$ws.Range("K141").Value = 6646.6158
$ws.Range("M141").Value = -1466.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9372.391
$ws.Range("I32").Value = 3410.6785
$ws.Range("K32").Value = 3410.6785
$ws.Range("M32").Value = -3123.6785
$ws.Range("H34").Value = 252500
$ws.Range("J34").Value = 252500
$ws.Range("L34").Value = 252500
$ws.Range("N34").Value = -253042
$ws.Range("H45").Value = 7815249
$ws.Range("I45").Value = 2168.5715
$ws.Range("J45").Value = 13892089
$ws.Range("K45").Value = 2168.5715
$ws.Range("L45").Value = 13892089
$ws.Range("M45").Value = -1791.5715
$ws.Range("N45").Value = -13892843
$ws.Range("H61").Value = 42319.56
$ws.Range("I61").Value = 2265.55
$ws.Range("K61").Value = 2265.55
$ws.Range("M61").Value = -2053.55
$ws.Range("H80").Value = 77450
$ws.Range("J80").Value = 77450
$ws.Range("L80").Value = 77450
$ws.Range("N80").Value = -79446
$ws.Range("H83").Value = 77450
$ws.Range("J83").Value = 77450
$ws.Range("L83").Value = 232350
$ws.Range("N83").Value = -242334
$ws.Range("H122").Value = 2630
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 3383.3333
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 10149.9999
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -15049.9999
$ws.Range("H128").Value = 68000
$ws.Range("J128").Value = 68000
$ws.Range("L128").Value = 68000
$ws.Range("N128").Value = -77960
$ws.Range("H136").Value = 42319.56
$ws.Range("I136").Value = 2265.55
$ws.Range("K136").Value = 6796.650000000001
$ws.Range("M136").Value = -4246.650000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1335.0769
$ws.Range("J58").Value = 1234.2
$ws.Range("L58").Value = 1234.2
$ws.Range("N58").Value = -1640.2
$ws.Range("H99").Value = 23729304
$ws.Range("I99").Value = 55558064
$ws.Range("J99").Value = 7814925
$ws.Range("K99").Value = 55558064
$ws.Range("L99").Value = 7814925
$ws.Range("M99").Value = -55556566
$ws.Range("N99").Value = -7817921
$ws.Range("H107").Value = 1482.44
$ws.Range("I107").Value = 1393.4762
$ws.Range("K107").Value = 1393.4762
$ws.Range("M107").Value = 526.5237999999999
$ws.Range("H122").Value = 3210.3157
$ws.Range("I122").Value = 2245.4443
$ws.Range("J122").Value = 4078.7
$ws.Range("K122").Value = 6736.3329
$ws.Range("L122").Value = 12236.1
$ws.Range("M122").Value = -4286.3329
$ws.Range("N122").Value = -17136.1
$ws.Range("H126").Value = 23729304
$ws.Range("I126").Value = 55558064
$ws.Range("J126").Value = 7814925
$ws.Range("K126").Value = 166674192
$ws.Range("L126").Value = 23444775
$ws.Range("M126").Value = -166671722
$ws.Range("N126").Value = -23449715
$ws.Range("H132").Value = 1950337.2
$ws.Range("I132").Value = 1895765.9
$ws.Range("J132").Value = 2168622.5
$ws.Range("K132").Value = 5687297.699999999
$ws.Range("L132").Value = 6505867.5
$ws.Range("M132").Value = -5684767.699999999
$ws.Range("N132").Value = -6510927.5
$ws.Range("H134").Value = 2899832
$ws.Range("I134").Value = 3762160.8
$ws.Range("K134").Value = 11286482.4
$ws.Range("M134").Value = -11283947.4
$ws.Range("H136").Value = 1335.0769
$ws.Range("J136").Value = 1234.2
$ws.Range("L136").Value = 3702.6
$ws.Range("N136").Value = -8802.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1537.125
$ws.Range("I102").Value = 1537.125
$ws.Range("K102").Value = 1537.125
$ws.Range("M102").Value = 84.875
$ws.Range("H122").Value = 5348762
$ws.Range("I122").Value = 6238917
$ws.Range("K122").Value = 18716751
$ws.Range("M122").Value = -18714301
$ws.Range("H126").Value = 3859.125
$ws.Range("I126").Value = 2576.75
$ws.Range("J126").Value = 5141.5
$ws.Range("K126").Value = 7730.25
$ws.Range("L126").Value = 15424.5
$ws.Range("M126").Value = -5260.25
$ws.Range("N126").Value = -20364.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4666.4736
$ws.Range("I7").Value = 4468
$ws.Range("J7").Value = 4782.25
$ws.Range("K7").Value = 4468
$ws.Range("L7").Value = 4782.25
$ws.Range("M7").Value = -4356
$ws.Range("N7").Value = -5006.25
$ws.Range("H16").Value = 1592.3226
$ws.Range("J16").Value = 1981.25
$ws.Range("L16").Value = 1981.25
$ws.Range("N16").Value = -2321.25
$ws.Range("H40").Value = 5053036.5
$ws.Range("J40").Value = 11113983
$ws.Range("L40").Value = 11113983
$ws.Range("N40").Value = -11114255
$ws.Range("H104").Value = 21444.5
$ws.Range("J104").Value = 21444.5
$ws.Range("L104").Value = 21444.5
$ws.Range("N104").Value = -28432.5
$ws.Range("H122").Value = 38713720
$ws.Range("I122").Value = 50003610
$ws.Range("J122").Value = 18186656
$ws.Range("K122").Value = 150010830
$ws.Range("L122").Value = 54559968
$ws.Range("M122").Value = -150008380
$ws.Range("N122").Value = -54564868
$ws.Range("H126").Value = 4666.4736
$ws.Range("I126").Value = 4468
$ws.Range("J126").Value = 4782.25
$ws.Range("K126").Value = 13404
$ws.Range("L126").Value = 14346.75
$ws.Range("M126").Value = -10934
$ws.Range("N126").Value = -19286.75
$ws.Range("H132").Value = 3159.6667
$ws.Range("I132").Value = 2540.625
$ws.Range("J132").Value = 3654.9
$ws.Range("K132").Value = 7621.875
$ws.Range("L132").Value = 10964.7
$ws.Range("M132").Value = -5091.875
$ws.Range("N132").Value = -16024.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1962.4375
$ws.Range("I122").Value = 1959
$ws.Range("J122").Value = 1964.5
$ws.Range("K122").Value = 5877
$ws.Range("L122").Value = 5893.5
$ws.Range("M122").Value = -3427
$ws.Range("N122").Value = -10793.5
$ws.Range("H124").Value = 50000
$ws.Range("J124").Value = 50000
$ws.Range("L124").Value = 50000
$ws.Range("N124").Value = -59820
$ws.Range("H132").Value = 2746.389
$ws.Range("I132").Value = 2536.4167
$ws.Range("J132").Value = 3166.3333
$ws.Range("K132").Value = 7609.250100000001
$ws.Range("L132").Value = 9498.999899999999
$ws.Range("M132").Value = -5079.250100000001
$ws.Range("N132").Value = -14558.9999
